$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits right after the mentor's
#    e-mail address (an artifact of whatever the last edit in the
#    *previous* save was). Word always re-stamps this bookmark at the
#    location of the most recent edit, so we pick it up from there and
#    will drop it back down at the spot we're about to edit below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

# ------------------------------------------------------------------
# 2) Fill in signature line 3) with the signer's name and the date,
#    underlining just the name and the date (matching lines 1/2's
#    blank style, and the already-signed entries elsewhere in the
#    document).
# ------------------------------------------------------------------
$target = "3) ___________________________________________________date__________________"

$findRange = $d.Content
$ok = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $ok) {
    throw "Could not find the blank signature line to fill in."
}

$lineStart = $findRange.Start

# Wipe the placeholder text; we'll rebuild it run-by-run so the name
# and date can carry their own (underlined) character formatting.
$findRange.Text = ""

$cursor = $lineStart

$partLabel = "3) ___"
$r = $d.Range($cursor, $cursor)
$r.InsertAfter($partLabel)
$cursor = $cursor + $partLabel.Length

$nameStart = $cursor
$partName = "Michele Tokuno"
$r = $d.Range($cursor, $cursor)
$r.InsertAfter($partName)
$cursor = $cursor + $partName.Length

$partMid = "__________________________________date___"
$r = $d.Range($cursor, $cursor)
$r.InsertAfter($partMid)
$cursor = $cursor + $partMid.Length

$dateStart = $cursor
$partDate = "03/31/2016"
$r = $d.Range($cursor, $cursor)
$r.InsertAfter($partDate)
$cursor = $cursor + $partDate.Length
$dateEnd = $cursor

$partTail = "_____"
$r = $d.Range($cursor, $cursor)
$r.InsertAfter($partTail)
$cursor = $cursor + $partTail.Length

# Underline just the signer's typed name and the typed date, after the
# fact, so the formatting doesn't leak onto the surrounding underscores.
$nameRange = $d.Range($nameStart, $nameStart + $partName.Length)
$nameRange.Font.Underline = 1

$dateRange = $d.Range($dateStart, $dateEnd)
$dateRange.Font.Underline = 1

# ------------------------------------------------------------------
# 3) Drop the "_GoBack" bookmark back in at the point of this edit
#    (right after the typed date), exactly like Word does whenever
#    you save after typing something.
# ------------------------------------------------------------------
$newGoBackRange = $d.Range($dateEnd, $dateEnd)
$d.Bookmarks.Add("_GoBack", $newGoBackRange)

Write-Output "Signed line 3) with 'Michele Tokuno' / 03/31/2016 and moved _GoBack bookmark."
